$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.455765
$ws.Range("H2").Value = 13.367295
$ws.Range("I2").Value = 0.1558824083674925
$ws.Range("J2").Value = 0.167793131187596
$ws.Range("M2").Value = 15.70818033333333
$ws.Range("N2").Value = 47.12454099999999
$ws.Range("O2").Value = 0.3220467100482788
$ws.Range("P2").Value = 0.334408980496766
$ws.Range("Q2").Value = 69.99196014295498
$ws.Range("R2").Value = 629.9276412865948
$ws.Range("S2").Value = 0.05020141676915325
$ws.Range("T2").Value = 0.05611152993480411

$ws.Range("G3").Value = 4.455765
$ws.Range("H3").Value = 13.367295
$ws.Range("I3").Value = 0.1558824083674925
$ws.Range("J3").Value = 0.167793131187596
$ws.Range("O3").Value = 0.1200026410479322
$ws.Range("P3").Value = 0.1246091315254933
$ws.Range("Q3").Value = 26.08075104389999
$ws.Range("R3").Value = 234.7267593951
$ws.Range("S3").Value = 0.01870630069701139
$ws.Range("T3").Value = 0.02090855635322951

$ws.Range("G4").Value = 4.455765
$ws.Range("H4").Value = 13.367295
$ws.Range("I4").Value = 0.1558824083674925
$ws.Range("J4").Value = 0.167793131187596
$ws.Range("M4").Value = 10.959131
$ws.Range("N4").Value = 32.877393
$ws.Range("O4").Value = 0.2246824271585863
$ws.Range("P4").Value = 0.2333072161810874
$ws.Range("Q4").Value = 48.831312340215
$ws.Range("R4").Value = 439.481811061935
$ws.Range("S4").Value = 0.03502403786333414
$ws.Range("T4").Value = 0.03914734833168602

$ws.Range("G5").Value = 4.455765
$ws.Range("H5").Value = 13.367295
$ws.Range("I5").Value = 0.1558824083674925
$ws.Range("J5").Value = 0.167793131187596
$ws.Range("M5").Value = 5.4093935
$ws.Range("N5").Value = 10.818787
$ws.Range("O5").Value = 0.1109025579706895
$ws.Range("P5").Value = 0.07677315161290731
$ws.Range("Q5").Value = 24.1029862285275
$ws.Range("R5").Value = 144.617917371165
$ws.Range("S5").Value = 0.01728775783058653
$ws.Range("T5").Value = 0.01288200750026976

$ws.Range("G6").Value = 4.455765
$ws.Range("H6").Value = 13.367295
$ws.Range("I6").Value = 0.1558824083674925
$ws.Range("J6").Value = 0.167793131187596
$ws.Range("M6").Value = 10.84612833333333
$ws.Range("N6").Value = 32.538385
$ws.Range("O6").Value = 0.2223656637745133
$ws.Range("P6").Value = 0.230901520183746
$ws.Range("Q6").Value = 48.32779901317499
$ws.Range("R6").Value = 434.9501911185749
$ws.Range("S6").Value = 0.03466289520740721
$ws.Range("T6").Value = 0.03874368906760665

$ws.Range("I7").Value = 0.1858758098371279
$ws.Range("J7").Value = 0.2000782799754709
$ws.Range("M7").Value = 15.70818033333333
$ws.Range("N7").Value = 47.12454099999999
$ws.Range("O7").Value = 0.3220467100482788
$ws.Range("P7").Value = 0.334408980496766
$ws.Range("Q7").Value = 83.45914340115355
$ws.Range("R7").Value = 751.1322906103819
$ws.Range("S7").Value = 0.05986069303560653
$ws.Range("T7").Value = 0.06690797362614374

$ws.Range("I8").Value = 0.1858758098371279
$ws.Range("J8").Value = 0.2000782799754709
$ws.Range("O8").Value = 0.1200026410479322
$ws.Range("P8").Value = 0.1246091315254933
$ws.Range("S8").Value = 0.02230558808737857
$ws.Range("T8").Value = 0.02493158070485793

$ws.Range("I9").Value = 0.1858758098371279
$ws.Range("J9").Value = 0.2000782799754709
$ws.Range("M9").Value = 10.959131
$ws.Range("N9").Value = 32.877393
$ws.Range("O9").Value = 0.2246824271585863
$ws.Range("P9").Value = 0.2333072161810874
$ws.Range("Q9").Value = 58.22696622218734
$ws.Range("R9").Value = 524.0426959996861
$ws.Range("S9").Value = 0.04176302810427373
$ws.Range("T9").Value = 0.04667970651937731

$ws.Range("I10").Value = 0.1858758098371279
$ws.Range("J10").Value = 0.2000782799754709
$ws.Range("M10").Value = 5.4093935
$ws.Range("N10").Value = 10.818787
$ws.Range("O10").Value = 0.1109025579706895
$ws.Range("P10").Value = 0.07677315161290731
$ws.Range("Q10").Value = 28.74065221111234
$ws.Range("R10").Value = 172.443913266674
$ws.Range("S10").Value = 0.02061410277581093
$ws.Range("T10").Value = 0.01536064012300654

$ws.Range("I11").Value = 0.1858758098371279
$ws.Range("J11").Value = 0.2000782799754709
$ws.Range("M11").Value = 10.84612833333333
$ws.Range("N11").Value = 32.538385
$ws.Range("O11").Value = 0.2223656637745133
$ws.Range("P11").Value = 0.230901520183746
$ws.Range("Q11").Value = 57.62657167858556
$ws.Range("R11").Value = 518.63914510727
$ws.Range("S11").Value = 0.04133239783405814
$ws.Range("T11").Value = 0.04619837900208537

$ws.Range("G12").Value = 5.633732333333334
$ws.Range("H12").Value = 16.901197
$ws.Range("I12").Value = 0.1970929266282699
$ws.Range("J12").Value = 0.2121524785267629
$ws.Range("M12").Value = 15.70818033333333
$ws.Range("N12").Value = 47.12454099999999
$ws.Range("O12").Value = 0.3220467100482788
$ws.Range("P12").Value = 0.334408980496766
$ws.Range("Q12").Value = 88.49568344173078
$ws.Range("R12").Value = 796.4611509755769
$ws.Range("S12").Value = 0.06347312859442113
$ws.Range("T12").Value = 0.07094569405399682

$ws.Range("G13").Value = 5.633732333333334
$ws.Range("H13").Value = 16.901197
$ws.Range("I13").Value = 0.1970929266282699
$ws.Range("J13").Value = 0.2121524785267629
$ws.Range("O13").Value = 0.1200026410479322
$ws.Range("P13").Value = 0.1246091315254933
$ws.Range("Q13").Value = 32.97570011740667
$ws.Range("R13").Value = 296.78130105666
$ws.Range("S13").Value = 0.02365167172725872
$ws.Range("T13").Value = 0.0264361361002008

$ws.Range("G14").Value = 5.633732333333334
$ws.Range("H14").Value = 16.901197
$ws.Range("I14").Value = 0.1970929266282699
$ws.Range("J14").Value = 0.2121524785267629
$ws.Range("M14").Value = 10.959131
$ws.Range("N14").Value = 32.877393
$ws.Range("O14").Value = 0.2246824271585863
$ws.Range("P14").Value = 0.2333072161810874
$ws.Range("Q14").Value = 61.74081065993568
$ws.Range("R14").Value = 555.667295939421
$ws.Range("S14").Value = 0.04428331713062886
$ws.Range("T14").Value = 0.04949670417099697

$ws.Range("G15").Value = 5.633732333333334
$ws.Range("H15").Value = 16.901197
$ws.Range("I15").Value = 0.1970929266282699
$ws.Range("J15").Value = 0.2121524785267629
$ws.Range("M15").Value = 5.4093935
$ws.Range("N15").Value = 10.818787
$ws.Range("O15").Value = 0.1109025579706895
$ws.Range("P15").Value = 0.07677315161290731
$ws.Range("Q15").Value = 30.47507506467317
$ws.Range("R15").Value = 182.850450388039
$ws.Range("S15").Value = 0.02185810972100455
$ws.Range("T15").Value = 0.01628761439898923

$ws.Range("G16").Value = 5.633732333333334
$ws.Range("H16").Value = 16.901197
$ws.Range("I16").Value = 0.1970929266282699
$ws.Range("J16").Value = 0.2121524785267629
$ws.Range("M16").Value = 10.84612833333333
$ws.Range("N16").Value = 32.538385
$ws.Range("O16").Value = 0.2223656637745133
$ws.Range("P16").Value = 0.230901520183746
$ws.Range("Q16").Value = 61.10418388298277
$ws.Range("R16").Value = 549.937654946845
$ws.Range("S16").Value = 0.04382669945495668
$ws.Range("T16").Value = 0.04898632980257908

$ws.Range("G17").Value = 6.087099
$ws.Range("H17").Value = 12.174198
$ws.Range("I17").Value = 0.2129537020222914
$ws.Range("J17").Value = 0.15281676675182
$ws.Range("M17").Value = 15.70818033333333
$ws.Range("N17").Value = 47.12454099999999
$ws.Range("O17").Value = 0.3220467100482788
$ws.Range("P17").Value = 0.334408980496766
$ws.Range("Q17").Value = 95.617248798853
$ws.Range("R17").Value = 573.7034927931179
$ws.Range("S17").Value = 0.06858103912888044
$ws.Range("T17").Value = 0.05110329917228822

$ws.Range("G18").Value = 6.087099
$ws.Range("H18").Value = 12.174198
$ws.Range("I18").Value = 0.2129537020222914
$ws.Range("J18").Value = 0.15281676675182
$ws.Range("O18").Value = 0.1200026410479322
$ws.Range("P18").Value = 0.1246091315254933
$ws.Range("Q18").Value = 35.62937309274
$ws.Range("R18").Value = 213.77623855644
$ws.Range("S18").Value = 0.02555500666360936
$ws.Range("T18").Value = 0.01904236458747818

$ws.Range("G19").Value = 6.087099
$ws.Range("H19").Value = 12.174198
$ws.Range("I19").Value = 0.2129537020222914
$ws.Range("J19").Value = 0.15281676675182
$ws.Range("M19").Value = 10.959131
$ws.Range("N19").Value = 32.877393
$ws.Range("O19").Value = 0.2246824271585863
$ws.Range("P19").Value = 0.2333072161810874
$ws.Range("Q19").Value = 66.70931535096901
$ws.Range("R19").Value = 400.2558921058141
$ws.Range("S19").Value = 0.04784695464277479
$ws.Range("T19").Value = 0.03565325443666167

$ws.Range("G20").Value = 6.087099
$ws.Range("H20").Value = 12.174198
$ws.Range("I20").Value = 0.2129537020222914
$ws.Range("J20").Value = 0.15281676675182
$ws.Range("M20").Value = 5.4093935
$ws.Range("N20").Value = 10.818787
$ws.Range("O20").Value = 0.1109025579706895
$ws.Range("P20").Value = 0.07677315161290731
$ws.Range("Q20").Value = 32.9275137644565
$ws.Range("R20").Value = 131.710055057826
$ws.Range("S20").Value = 0.02361711028360011
$ws.Range("T20").Value = 0.01173222480283177

$ws.Range("G21").Value = 6.087099
$ws.Range("H21").Value = 12.174198
$ws.Range("I21").Value = 0.2129537020222914
$ws.Range("J21").Value = 0.15281676675182
$ws.Range("M21").Value = 10.84612833333333
$ws.Range("N21").Value = 32.538385
$ws.Range("O21").Value = 0.2223656637745133
$ws.Range("P21").Value = 0.230901520183746
$ws.Range("Q21").Value = 66.021456931705
$ws.Range("R21").Value = 396.12874159023
$ws.Range("S21").Value = 0.04735359130342674
$ws.Range("T21").Value = 0.03528562375256017

$ws.Range("G22").Value = 7.094445666666666
$ws.Range("H22").Value = 21.283337
$ws.Range("I22").Value = 0.2481951531448182
$ws.Range("J22").Value = 0.2671593435583502
$ws.Range("M22").Value = 15.70818033333333
$ws.Range("N22").Value = 47.12454099999999
$ws.Range("O22").Value = 0.3220467100482788
$ws.Range("P22").Value = 0.334408980496766
$ws.Range("Q22").Value = 111.4408318970352
$ws.Range("R22").Value = 1002.967487073317
$ws.Range("S22").Value = 0.07993043252021743
$ws.Range("T22").Value = 0.08934048370953315

$ws.Range("G23").Value = 7.094445666666666
$ws.Range("H23").Value = 21.283337
$ws.Range("I23").Value = 0.2481951531448182
$ws.Range("J23").Value = 0.2671593435583502
$ws.Range("O23").Value = 0.1200026410479322
$ws.Range("P23").Value = 0.1246091315254933
$ws.Range("Q23").Value = 41.52563504287333
$ws.Range("R23").Value = 373.73071538586
$ws.Range("S23").Value = 0.02978407387267419
$ws.Range("T23").Value = 0.03329049377972692

$ws.Range("G24").Value = 7.094445666666666
$ws.Range("H24").Value = 21.283337
$ws.Range("I24").Value = 0.2481951531448182
$ws.Range("J24").Value = 0.2671593435583502
$ws.Range("M24").Value = 10.959131
$ws.Range("N24").Value = 32.877393
$ws.Range("O24").Value = 0.2246824271585863
$ws.Range("P24").Value = 0.2333072161810874
$ws.Range("Q24").Value = 77.74895943338234
$ws.Range("R24").Value = 699.7406349004411
$ws.Range("S24").Value = 0.0557650894175748
$ws.Range("T24").Value = 0.0623302027223654

$ws.Range("G25").Value = 7.094445666666666
$ws.Range("H25").Value = 21.283337
$ws.Range("I25").Value = 0.2481951531448182
$ws.Range("J25").Value = 0.2671593435583502
$ws.Range("M25").Value = 5.4093935
$ws.Range("N25").Value = 10.818787
$ws.Range("O25").Value = 0.1109025579706895
$ws.Range("P25").Value = 0.07677315161290731
$ws.Range("Q25").Value = 38.37664827536983
$ws.Range("R25").Value = 230.259889652219
$ws.Range("S25").Value = 0.02752547735968736
$ws.Range("T25").Value = 0.02051066478781001

$ws.Range("G26").Value = 7.094445666666666
$ws.Range("H26").Value = 21.283337
$ws.Range("I26").Value = 0.2481951531448182
$ws.Range("J26").Value = 0.2671593435583502
$ws.Range("M26").Value = 10.84612833333333
$ws.Range("N26").Value = 32.538385
$ws.Range("O26").Value = 0.2223656637745133
$ws.Range("P26").Value = 0.230901520183746
$ws.Range("Q26").Value = 76.94726815452721
$ws.Range("R26").Value = 692.5254133907449
$ws.Range("S26").Value = 0.05519007997466448
$ws.Range("T26").Value = 0.06168749855891473
